$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common values shared by every new row (72-77)
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$fecha     = 44476
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100107
$producto  = "Otros"
$categoriaId = 100107002
$categoria = "Chirimoya"
$variedad  = "Cultivar IV Región"
$origen    = "Provincia de Limarí"

# Per-row data: Calidad, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado,
#               Unidad de comercializacion, Precio $/Kg, Kg / unidad
$rows = @(
    @{ Row=72; Calidad="Cuarta";                  Volumen=350; PMin=1000;  PMax=1000;  PProm=1000;  Unidad="`$/kilo (en caja de 15 kilos)"; PrecioKg=1000; KgUnidad=1 },
    @{ Row=73; Calidad="Especial";                Volumen=280; PMin=20000; PMax=20000; PProm=20000; Unidad="`$/bandeja 8 kilos";              PrecioKg=2500; KgUnidad=8 },
    @{ Row=74; Calidad="Extra (doble especial)";   Volumen=220; PMin=24000; PMax=24000; PProm=24000; Unidad="`$/bandeja 8 kilos";              PrecioKg=3000; KgUnidad=8 },
    @{ Row=75; Calidad="Primera";                  Volumen=300; PMin=16000; PMax=16000; PProm=16000; Unidad="`$/bandeja 8 kilos";              PrecioKg=2000; KgUnidad=8 },
    @{ Row=76; Calidad="Segunda";                  Volumen=300; PMin=12800; PMax=12800; PProm=12800; Unidad="`$/bandeja 8 kilos";              PrecioKg=1600; KgUnidad=8 },
    @{ Row=77; Calidad="Tercera";                  Volumen=310; PMin=1400;  PMax=1400;  PProm=1400;  Unidad="`$/kilo (en caja de 15 kilos)";   PrecioKg=1400; KgUnidad=1 }
)

foreach ($r in $rows) {
    $i = $r.Row

    $ws.Cells.Item($i, 1).Value  = $mercadoId
    $ws.Cells.Item($i, 2).Value  = $mercado
    $ws.Cells.Item($i, 3).Value  = $region

    $ws.Cells.Item($i, 4).Value  = $fecha
    $ws.Cells.Item($i, 4).NumberFormat = $ws.Cells.Item($i - 1, 4).NumberFormat

    $ws.Cells.Item($i, 5).Value  = $codreg
    $ws.Cells.Item($i, 6).Value  = $tipo
    $ws.Cells.Item($i, 7).Value  = $productoId
    $ws.Cells.Item($i, 8).Value  = $producto
    $ws.Cells.Item($i, 9).Value  = $categoriaId
    $ws.Cells.Item($i, 10).Value = $categoria
    $ws.Cells.Item($i, 11).Value = $variedad
    $ws.Cells.Item($i, 12).Value = $r.Calidad
    $ws.Cells.Item($i, 13).Value = $r.Volumen
    $ws.Cells.Item($i, 14).Value = $r.PMin
    $ws.Cells.Item($i, 15).Value = $r.PMax
    $ws.Cells.Item($i, 16).Value = $r.PProm
    $ws.Cells.Item($i, 17).Value = $r.Unidad
    $ws.Cells.Item($i, 18).Value = $origen
    $ws.Cells.Item($i, 19).Value = $r.PrecioKg
    $ws.Cells.Item($i, 20).Value = $r.KgUnidad
}
